$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "Laura Sofia Baron Molina" (row 4),
# shifting the rows below it up by one.
$ws.Rows.Item(4).Delete()
